$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, pushing the existing weekly records
# (rows 14-27) down to rows 15-28 while keeping their data intact.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with this week's price record.
$ws.Cells.Item(14, 1).Value = 4
$ws.Cells.Item(14, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(14, 3).Value = "Los Lagos"
$ws.Cells.Item(14, 4).Value = 44810
$ws.Cells.Item(14, 5).Value = 10
$ws.Cells.Item(14, 6).Value = 100112013
$ws.Cells.Item(14, 7).Value = "Alcachofa"
$ws.Cells.Item(14, 8).Value = "Madrigal"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 100
$ws.Cells.Item(14, 11).Value = 13500
$ws.Cells.Item(14, 12).Value = 14000
$ws.Cells.Item(14, 13).Value = 13750
$ws.Cells.Item(14, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(14, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(14, 16).Value = 344
$ws.Cells.Item(14, 17).Value = 40
$ws.Cells.Item(14, 18).Value = "Hortaliza"
